$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update electricity_rate (col G) values
$ws.Range("G19").Value = 55.2
$ws.Range("G24").Value = 57
$ws.Range("G25").Value = 54.4
$ws.Range("G26").Value = 54.4

# Replace set_voltage (col H) formulas with hardcoded values for rows 25 and 26
$ws.Range("H25").Value = 53
$ws.Range("H26").Value = 52

# Update the selected/active cell shown in the sheet view
$ws.Range("G27").Select()
